# Update cryptocurrency price/volume snapshot cells to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value is a plain decimal number need to be force-typed
# as Text (matching the rest of the Price column, which stores values like
# "37.742.95" / "0.0₃0834" as literal strings, not numbers) so Excel does not
# silently reinterpret them as numeric values.
$textCells = @(
    "D5",
    "D8",
    "D9",
    "D13",
    "D14",
    "D15",
    "D22",
    "D26",
    "D28",
    "D29",
    "D35",
    "D38",
    "D39",
    "D40",
    "D41",
    "D47",
    "D49"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '37.749.35'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = '2.077.19'
$ws.Range("E3").Value = '  -1.23%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '233.60'
$ws.Range("E5").Value = '  -0.66%  '

$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '58.28'
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").Value = '0.395'
$ws.Range("E9").Value = '  +1.07%  '

$ws.Range("E10").Value = '  +0.84%  '

$ws.Range("D12").Value = '2.382.26'
$ws.Range("E12").Value = '  -1.30%  '

$ws.Range("D13").Value = '14.84'
$ws.Range("E13").Value = '  +2.50%  '

$ws.Range("D14").Value = '20.87'
$ws.Range("E14").Value = '  -1.94%  '

$ws.Range("D15").Value = '0.773'
$ws.Range("E15").Value = '  -1.40%  '

$ws.Range("E16").Value = '  +1.51%  '

$ws.Range("D17").Value = '2.085.72'
$ws.Range("E17").Value = '  -0.73%  '

$ws.Range("D18").Value = '37.674.70'
$ws.Range("E18").Value = '  -0.30%  '

$ws.Range("E19").Value = '  -0.91%  '

$ws.Range("E20").Value = '  +1.21%  '

$ws.Range("D21").Value = '0.0₃0835'
$ws.Range("E21").Value = '  +1.49%  '

$ws.Range("D22").Value = '228.54'
$ws.Range("E22").Value = '  +0.50%  '

$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("E24").Value = '  -0.64%  '

$ws.Range("E25").Value = '  -0.46%  '

$ws.Range("D26").Value = '170.40'
$ws.Range("E26").Value = '  +1.52%  '

$ws.Range("E27").Value = '  +4.23%  '

$ws.Range("D28").Value = '9.03'
$ws.Range("E28").Value = '  +0.88%  '

$ws.Range("D29").Value = '19.46'
$ws.Range("E29").Value = '  -0.43%  '

$ws.Range("E30").Value = '  -2.10%  '

$ws.Range("E31").Value = '  +2.64%  '

$ws.Range("E32").Value = '  +0.79%  '

$ws.Range("E33").Value = '  +1.20%  '

$ws.Range("E34").Value = '  +1.22%  '

$ws.Range("D35").Value = '2.48'
$ws.Range("E35").Value = '  -4.69%  '

$ws.Range("E36").Value = '  +2.52%  '

$ws.Range("E37").Value = '  -3.40%  '

$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.07%  '

$ws.Range("D39").Value = '5.28'
$ws.Range("E39").Value = '  -2.79%  '

$ws.Range("D40").Value = '0.0972'
$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("D41").Value = '98.08'
$ws.Range("E41").Value = '  +1.48%  '

$ws.Range("E42").Value = '  -2.13%  '

$ws.Range("E43").Value = '  +0.54%  '

$ws.Range("D44").Value = '1.450.21'
$ws.Range("E44").Value = '  -1.51%  '

$ws.Range("E45").Value = '  +2.46%  '

$ws.Range("E46").Value = '  -1.14%  '

$ws.Range("D47").Value = '16.32'
$ws.Range("E47").Value = '  +5.82%  '

$ws.Range("E48").Value = '  +0.78%  '

$ws.Range("D49").Value = '7.40'
$ws.Range("E49").Value = '  +1.04%  '

$ws.Range("E50").Value = '  -0.53%  '

$ws.Range("D51").Value = '2.266.66'
